$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.210.28"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.439.06"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  +1.68%  "
$ws.Range("D5").Value = "'0.9139"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.00%  "
$ws.Range("D6").Value = "'275.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").Value = "'0.3620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").Value = "'0.3077"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'1.024"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").Value = "'0.06491"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "'0.9993"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'5.337"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").Value = "'17.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'6.044"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "1.438.80"
$ws.Range("E17").Value = "  +2.48%  "
$ws.Range("D18").Value = "'0.9329"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.31%  "
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "'67.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.54%  "
$ws.Range("D21").Value = "'5.405"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("D22").Value = "'14.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.71%  "
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "'2.241"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "20.222.14"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("D26").Value = "'138.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "'2.125"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.20%  "
$ws.Range("D28").Value = "'16.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").Value = "1.590.14"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").Value = "'110.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "'3.833"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.56%  "
$ws.Range("D32").Value = "'0.8051"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("D33").Value = "'4.825"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.77%  "
$ws.Range("D34").Value = "'0.07642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "'1.475"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").Value = "'0.05844"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "'4.653"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").Value = "'1.128"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.32%  "
$ws.Range("D39").Value = "'0.01984"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.03%  "
$ws.Range("D40").Value = "'10.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("D41").Value = "'0.1848"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.59%  "
$ws.Range("D42").Value = "'0.9280"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.76%  "
$ws.Range("D43").Value = "'7.179"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -14.71%  "
$ws.Range("D44").Value = "'0.5203"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D45").Value = "'3.487"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "'11.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").Value = "'116.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.66%  "
$ws.Range("D48").Value = "'0.5076"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "'1.731"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("D50").Value = "'0.06338"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("E51").Value = "  -0.97%  "
